# Optuna Attempt (go back with original)
# Update forecast values on the "Forecast Comparison" sheet and the
# corresponding rolled-up totals on the "Summary" sheet.

$wb = $excel.ActiveWorkbook
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet ---

# Row 2 (W8)
$wsForecast.Range("L2").Value = 0.89

# Row 3 (W9)
$wsForecast.Range("L3").Value = 1.11

# Row 4 (W10)
$wsForecast.Range("D4").Value = 73
$wsForecast.Range("H4").Value = 8.359999999999999
$wsForecast.Range("L4").Value = 1.05

# Row 5 (W11)
$wsForecast.Range("H5").Value = 4.05
$wsForecast.Range("L5").Value = 0.84

# Row 6 (W12)
$wsForecast.Range("H6").Value = 3.02
$wsForecast.Range("L6").Value = 0.97

# Row 7 (W13)
$wsForecast.Range("H7").Value = 2.06
$wsForecast.Range("L7").Value = 0.96

# Row 8 (W14)
$wsForecast.Range("H8").Value = 1.06
$wsForecast.Range("J8").Value = "Normal"
$wsForecast.Range("L8").Value = 1.13

# Row 9 (W15)
$wsForecast.Range("H9").Value = 0.06
$wsForecast.Range("L9").Value = 0.95

# Row 10 (W16)
$wsForecast.Range("L10").Value = 1.06

# Row 11 (W17)
$wsForecast.Range("L11").Value = 0.9399999999999999

# Row 12 (W18)
$wsForecast.Range("D12").Value = 66
$wsForecast.Range("L12").Value = 0.88

# Row 13 (W19)
$wsForecast.Range("D13").Value = 66
$wsForecast.Range("L13").Value = 1.18

# Row 14 (W20)
$wsForecast.Range("D14").Value = 65
$wsForecast.Range("L14").Value = 1

# Row 15 (W21)
$wsForecast.Range("D15").Value = 62
$wsForecast.Range("L15").Value = 1.13

# Row 16 (W22)
$wsForecast.Range("D16").Value = 63
$wsForecast.Range("L16").Value = 0.83

# Row 17 (W23)
$wsForecast.Range("D17").Value = 61
$wsForecast.Range("L17").Value = 1.14

# --- Summary sheet ---
# These cells hold numeric-looking values but are stored as text in the
# workbook, so force a text number format before assigning the string.

$wsSummary.Range("B9").NumberFormat = "@"
$wsSummary.Range("B9").Value = "1756"

$wsSummary.Range("B10").NumberFormat = "@"
$wsSummary.Range("B10").Value = "1117"

$wsSummary.Range("B11").NumberFormat = "@"
$wsSummary.Range("B11").Value = "586"

$wsSummary.Range("B14").NumberFormat = "@"
$wsSummary.Range("B14").Value = "62"
